$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add owner "Wim" and status "to be started"
$ws.Range("C7").Value = "Wim"
$ws.Range("D7").Value = "to be started"

# Fix typo in row 11 task description: "an feature matrix" -> "a feature matrix"
$ws.Range("B11").Value = "write a parsing script that takes the cluster buster output and produces a feature matrix`n"

# Update the active selection to D11 (as recorded in the saved view state)
$ws.Range("D11").Select()
